$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Taxon column (D): wrap the scientific (italicised) name in asterisks
$ws.Range("D2").Value  = "Black rhinoceros (*Diceros bicornis*)"
$ws.Range("D3").Value  = "African wild dog (*Lycaon pictus*)"
$ws.Range("D4").Value  = "Cheetah (*Acinonyx jubatus*)"
$ws.Range("D5").Value  = "Riverine rabbit (*Bunolagus monticularis*)"
$ws.Range("D7").Value  = "Black-footed cat (*Felis nigripes*)"
$ws.Range("D8").Value  = "Cape Mountain Zebra (*Equus zebra zebra*)"
$ws.Range("D9").Value  = "Bontebok (*Damaliscus pygargus*)"
$ws.Range("D10").Value = "Leopard (*Panthera pardus*)"

# Contact column (E): turn into Markdown-style links
$ws.Range("E2").Value  = "[Conservation: BRREP](wwf.org.za)"
$ws.Range("E3").Value  = " [Carnivore Conservation Programme](ewt.org)"
$ws.Range("E4").Value  = "[Carnivore Conservation Programme](ewt.org)"
$ws.Range("E5").Value  = "[ewt.org](ewt.org)"
$ws.Range("E6").Value  = "[pangolin.org.za](pangolin.org.za)"
$ws.Range("E7").Value  = "[black-footed-cat.wild-cat.org](black-footed-cat.wild-cat.org)"
$ws.Range("E8").Value  = "[capenature.co.za](capenature.co.za); [sanparks.org](sanparks.org)"
$ws.Range("E9").Value  = "[capenature.co.za](capenature.co.za); [sanparks.org](sanparks.org)"
$ws.Range("E10").Value = "[capeleopard.org.za](capeleopard.org.za)"

# Update the view: new active cell/selection (cursor moved from D11 to E10)
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E10").Select()
